$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Catalogo")
$ws2 = $wb.Worksheets.Item("SC")

# 1. Re-filter the Catalogo station list from state "RS" to state "PR".
#    Using the xlFilterValues-style overload (Operator = 7) reproduces the
#    same serialized <filters><filter val="PR"/></filters> markup Excel
#    emits for a simple value-list AutoFilter, and re-derives every row's
#    hidden/visible state exactly like a live filter change would.
$ws1.Range("A1:H572").AutoFilter(2, @("PR"), 7)

# 2. Select the whole row for the station that is about to be copied
#    (mirrors right-clicking/selecting row 209 before copying it).
$ws1.Range("A209:XFD209").Select()

# 3. Copy that station's row into the SC sheet as the new row 25 - formats
#    first (so the date cell reuses the existing short-date style instead
#    of registering a new number format), then values.
$ws1.Range("A209:H209").Copy()
$ws2.Range("A25:H25").PasteSpecial(-4122)
$ws1.Range("A209:H209").Copy()
$ws2.Range("A25:H25").PasteSpecial(-4163)

# 4. Leave the workbook focused on the SC sheet with the same selection
#    the author ended up on.
$ws2.Activate()
$ws2.Range("E31").Select()
